# Risk Table.xlsx - 2nd semester documentation update.
# Remove the "Usability Study" risk (website is not user friendly) and the
# "server risk" (cannot get siue email / host server) rows from the
# active risk table on the "Sp2" worksheet. All rows below shift up and
# the table / dimension / conditional-formatting ranges shrink to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sp2")
$ws.Activate()

# Original row 6 ("cannot get siue email for notification service" / server
# risk) and original row 3 ("website is not user friendly" / usability
# study) are being removed. Delete bottom-up so row indices of the
# not-yet-deleted row stay valid.
$ws.Range("A6").EntireRow.Delete()
$ws.Range("A3").EntireRow.Delete()

# The color-scale conditional formatting over column E was anchored to
# E1:E10; after removing two rows it should cover E1:E8.
$newDataRange = $ws.Range("E1:E8")
$conds = $newDataRange.FormatConditions
for ($i = 1; $i -le $conds.Count; $i++) {
    $conds.Item($i).ModifyAppliesToRange($newDataRange)
}

# Leave the selection on the (now) row 3, matching the last row-delete
# operation performed interactively.
$ws.Rows("3:3").Select() | Out-Null
